$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a plain-text value into a cell without Excel auto-converting
# it to a number/date (needed for columns that store numeric-looking text,
# e.g. "6.76", "0.0392", quarter labels, fund codes, ...).
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$val) {
    $range.NumberFormat = "@"
    $range.Value = $val
}

# ===========================================================================
# 1) Duplicate the "2022-Q3" sheet (current position 2) to create the new
#    "2022-Q4" sheet right before it, so it inherits identical headers,
#    column widths / cell styles and the A-column row-index formatting.
# ===========================================================================
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3) | Out-Null
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Overwrite the fund rows (2-5) on the new "2022-Q4" sheet with this
# quarter's holdings.
Set-TextValue $q4.Range("B2") "015784"
Set-TextValue $q4.Range("C2") "中信建投中证1000指数增强A"
Set-TextValue $q4.Range("D2") "6.76"
Set-TextValue $q4.Range("E2") "89.78"
Set-TextValue $q4.Range("F2") "0.58"
Set-TextValue $q4.Range("G2") "0.0392"
$q4.Range("H2").Value = 10

Set-TextValue $q4.Range("B3") "015785"
Set-TextValue $q4.Range("C3") "中信建投中证1000指数增强C"
Set-TextValue $q4.Range("D3") "2.40"
Set-TextValue $q4.Range("E3") "89.78"
Set-TextValue $q4.Range("F3") "0.58"
Set-TextValue $q4.Range("G3") "0.0139"
$q4.Range("H3").Value = 10

Set-TextValue $q4.Range("B4") "009327"
Set-TextValue $q4.Range("C4") "东兴兴晟混合A"
Set-TextValue $q4.Range("D4") "0.38"
Set-TextValue $q4.Range("E4") "79.79"
Set-TextValue $q4.Range("F4") "1.20"
Set-TextValue $q4.Range("G4") "0.0046"
$q4.Range("H4").Value = 1

Set-TextValue $q4.Range("B5") "009328"
Set-TextValue $q4.Range("C5") "东兴兴晟混合C"
Set-TextValue $q4.Range("D5") "0.08"
Set-TextValue $q4.Range("E5") "79.79"
Set-TextValue $q4.Range("F5") "1.20"
Set-TextValue $q4.Range("G5") "0.0010"
$q4.Range("H5").Value = 1

# The "@" number-format stamp used by Set-TextValue leaves a stray cell
# style behind (B:G originally carried no explicit style). Re-flatten it by
# pasting the (style-only) format of the still-untouched H column back onto
# B:G for each data row.
for ($r = 2; $r -le 5; $r++) {
    $q4.Range("H$r").Copy() | Out-Null
    $q4.Range("B$r`:G$r").PasteSpecial(-4122) | Out-Null
}

# ===========================================================================
# 2) Update the "总计" (summary) sheet: insert the new 2022-Q4 total as the
#    first data row, pushing the existing rows down by one.
# ===========================================================================
$total = $wb.Worksheets.Item("总计")

# Extend the existing row-5 formatting down to the new row 6 first (format
# only - values are written explicitly below).
$total.Range("A5:D5").Copy() | Out-Null
$total.Range("A6:D6").PasteSpecial(-4122) | Out-Null

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.06

$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2022-Q3"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
Set-TextValue $total.Range("B4") "2021-Q4"
$total.Range("C4").Value = 9
$total.Range("D4").Value = 0.97

$total.Range("A5").Value = 3
Set-TextValue $total.Range("B5") "2021-Q1"
$total.Range("C5").Value = 10
$total.Range("D5").Value = 1.61

$total.Range("A6").Value = 4
Set-TextValue $total.Range("B6") "2020-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.01

# Re-flatten the B-column style pollution introduced by the "@" stamp above
# (column C never had NumberFormat touched, so its style is the clean
# baseline to copy back onto B).
for ($r = 2; $r -le 6; $r++) {
    $total.Range("C$r").Copy() | Out-Null
    $total.Range("B$r").PasteSpecial(-4122) | Out-Null
}
